$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rating amount (D) / rating value (E) updates
$ws.Range("D5").Value = 22

$ws.Range("D16").Value = 13

$ws.Range("D33").Value = 10
$ws.Range("E33").Value = 4.5

$ws.Range("D34").Value = 6

$ws.Range("D46").Value = 7
$ws.Range("E46").Value = 4.5

$ws.Range("D47").Value = 14

$ws.Range("D48").Value = 21

$ws.Range("D51").Value = 7
$ws.Range("E51").Value = 4.5

$ws.Range("D60").Value = 3

$ws.Range("D66").Value = 10
$ws.Range("E66").Value = 4.5

$ws.Range("D114").Value = 7

$ws.Range("D130").Value = 14

# Product aria-label text updates ("- Online kein Bestand" inserted)
$ws.Range("M380").Value = "Leisi Kuchenteig rund ausgewallt Ø32cm glutenfrei - Online kein Bestand 4.95 Schweizer Franken"

$ws.Range("M384").Value = "Leisi Blätterteig glutenfrei rund ausgewallt Ø32cm - Online kein Bestand 4.95 Schweizer Franken"

$ws.Range("M386").Value = "Buitoni Pizzateig Glutenfrei &amp; Ohne Lactose Rund Ausgewallt Ø25cm - Online kein Bestand 4.95 Schweizer Franken"

# Refresh crawl timestamp across every data row
$ws.Range("O2:O399").Value = "2023-01-05 12:56:40"
